$d = $word.ActiveDocument

# The document has a single section whose header/footer stories each
# carry one inline picture (the Pearson logo in the footers, the BTec
# logo in the headers). Both the "primary" and "first page" flavours
# exist (indices 1 and 2), and each needs the same rename applied to
# the embedded picture's internal OOXML name
# (wp:docPr/@name + pic:cNvPr/@name) - the picture's description/alt
# text is left untouched, only the internal file-name-style label
# changes.

for ($i = 1; $i -le 2; $i++) {
    $ftr = $d.Sections(1).Footers($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

for ($i = 1; $i -le 2; $i++) {
    $hdr = $d.Sections(1).Headers($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
